$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.923.48'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.879.86'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3871'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07852'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9861'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.77'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('D12').Value = '1.860.78'
$ws.Range('E12').Value = '  -1.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.989'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.649'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06978'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.13'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.004'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009957'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('D21').Value = '28.912.53'
$ws.Range('E21').Value = '  +0.40%  '
$ws.Range('E22').Value = '  -1.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.101'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.85%  '
$ws.Range('E26').Value = '  -2.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.988'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.40%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '117.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.921'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09371'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.9012'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.266'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.318'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.255'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('E35').Value = '  +1.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05745'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02072'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.647'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5646'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1766'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.703'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.276'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5346'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07047'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.844'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.538'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '112.84'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.070'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.70'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.76%  '
